# Add new conversation flows (problems / stressed) to the SPGeTTi message
# workbook, plus the corresponding Follow Up rows that back them.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("User_Initiated_Messages")
$ws4 = $wb.Worksheets.Item("Follow_Up_Messages")

# ---------------------------------------------------------------------
# User_Initiated_Messages: two new top-level conversation starters
# ---------------------------------------------------------------------

# Row 3 - "problems" flow
$ws3.Cells.Item(3,1).Value = 1
$ws3.Cells.Item(3,2).Value = "problems"
$ws3.Cells.Item(3,3).Value = "Message;Buttons"
$ws3.Cells.Item(3,4).Value = "Do you think you have a problem?"
$ws3.Cells.Item(3,5).Value = "Yes;No"
$ws3.Cells.Item(3,6).Value = "8;9"

# ---------------------------------------------------------------------
# Follow_Up_Messages: new rows 10-20 supporting the "problems" flow
# ---------------------------------------------------------------------

# Row 10 (ID 8)
$ws4.Cells.Item(10,1).Value = 8
$ws4.Cells.Item(10,2).Value = "yes"
$ws4.Cells.Item(10,3).Value = "Message;Buttons"
$ws4.Cells.Item(10,4).Value = "You are not alone. Do you want to talk to someone?"
$ws4.Cells.Item(10,5).Value = "Yes;No"
$ws4.Cells.Item(10,6).Value = "10;11"

# Row 12 (ID 10)
$ws4.Cells.Item(12,1).Value = 10
$ws4.Cells.Item(12,2).Value = "yes"
$ws4.Cells.Item(12,3).Value = "Message"
$ws4.Cells.Item(12,4).Value = "You can phone a friend, call the Gambling Helpline 0800 654 655 or text 8006 for more information"

# Row 13 (ID 11)
$ws4.Cells.Item(13,1).Value = 11
$ws4.Cells.Item(13,2).Value = "no"
$ws4.Cells.Item(13,3).Value = "Message"
$ws4.Cells.Item(13,4).Value = "If you want more information, you can go to choicenotchance.org.nz"

# Row 11 (ID 9)
$ws4.Cells.Item(11,1).Value = 9
$ws4.Cells.Item(11,2).Value = "no"
$ws4.Cells.Item(11,3).Value = "Message;Buttons"
$d11 = $ws4.Cells.Item(11,4)
$d11.Value = "Do you think you can win back what you have lost by gambling more?"
$d11.WrapText = $true
$ws4.Cells.Item(11,5).Value = "Yes;No"
$ws4.Cells.Item(11,6).Value = 12

# Row 14 (ID 12)
$ws4.Cells.Item(14,1).Value = 12
$ws4.Cells.Item(14,3).Value = "Message;Buttons"
$ws4.Cells.Item(14,4).Value = "Do you gamble using money meant for other things?"
$ws4.Cells.Item(14,5).Value = "Yes;No"
$ws4.Cells.Item(14,6).Value = 13

# Row 15 (ID 13)
$ws4.Cells.Item(15,1).Value = 13
$ws4.Cells.Item(15,3).Value = "Message;Buttons"
$ws4.Cells.Item(15,4).Value = "Do you lie about your gambling to people you care about?"
$ws4.Cells.Item(15,5).Value = "Yes;No"
$ws4.Cells.Item(15,6).Value = 14

# Row 16 (ID 14)
$ws4.Cells.Item(16,1).Value = 14
$ws4.Cells.Item(16,3).Value = "Message;Buttons"
$ws4.Cells.Item(16,4).Value = "Do you get angry or defensive when people ask about your gambling?"
$ws4.Cells.Item(16,5).Value = "Yes;No"
$ws4.Cells.Item(16,6).Value = 15

# Row 17 (ID 15)
$ws4.Cells.Item(17,1).Value = 15
$ws4.Cells.Item(17,3).Value = "Message;Buttons"
$ws4.Cells.Item(17,4).Value = "Do you borrow money just so you can gamble?"
$ws4.Cells.Item(17,5).Value = "Yes;No"
$ws4.Cells.Item(17,6).Value = 16

# Row 18 (ID 16)
$ws4.Cells.Item(18,1).Value = 16
$ws4.Cells.Item(18,3).Value = "Message;Buttons"
$ws4.Cells.Item(18,4).Value = "Did you answer yes to most of these?"
$ws4.Cells.Item(18,5).Value = "Yes;No"
$ws4.Cells.Item(18,6).Value = "17;18"

# Row 19 (ID 17)
$ws4.Cells.Item(19,1).Value = 17
$ws4.Cells.Item(19,2).Value = "yes"
$ws4.Cells.Item(19,3).Value = "Message;Buttons"
$ws4.Cells.Item(19,4).Value = "It looks like gambling is hurting you and maybe those you care about. Do you want to talk to someone about it?"
$ws4.Cells.Item(19,5).Value = "Yes;No"
$ws4.Cells.Item(19,6).Value = "10;11"

# Row 20 (ID 18)
$ws4.Cells.Item(20,1).Value = 18
$ws4.Cells.Item(20,2).Value = "no"
$ws4.Cells.Item(20,3).Value = "Message"
$ws4.Cells.Item(20,4).Value = "Okay that's great! If you want more information, you can go to choicenotchance.org.nz"

# Clear the stale Follow Ups value on the "echo" lapse message - it no longer
# points anywhere now that the follow-up list has been extended.
$ws4.Cells.Item(7,6).ClearContents()

# ---------------------------------------------------------------------
# User_Initiated_Messages: "stressed" flow
# ---------------------------------------------------------------------

$ws3.Cells.Item(4,1).Value = 2
$ws3.Cells.Item(4,2).Value = "stressed"
$ws3.Cells.Item(4,3).Value = "Message;Buttons"
$ws3.Cells.Item(4,4).Value = "Do you need some ideas to help you relax?"
$ws3.Cells.Item(4,5).Value = "Yes;No"
$ws3.Cells.Item(4,6).Value = "19;20"

# ---------------------------------------------------------------------
# Follow_Up_Messages: new rows 21-24 supporting the "stressed" flow
# ---------------------------------------------------------------------

# Row 21 (ID 19)
$ws4.Cells.Item(21,1).Value = 19
$ws4.Cells.Item(21,2).Value = "yes"
$ws4.Cells.Item(21,3).Value = "Message;Buttons"
$ws4.Cells.Item(21,5).Value = "Indoors;Outdoors"
$ws4.Cells.Item(21,4).Value = "Do you prefer indoors or outdoors activities?"
$ws4.Cells.Item(21,6).Value = "21;22"

# Row 23/24 keys entered first ...
$ws4.Cells.Item(23,2).Value = "indoors"
$ws4.Cells.Item(24,2).Value = "outdoors"

# ... then their messages ...
$ws4.Cells.Item(23,4).Value = "Why don't you invite your friends over for dinner? You find some recipes at myfamily.kiwi/foods"
$ws4.Cells.Item(24,4).Value = "You can go for a walk, visit a museum or gallery with your friends or family. Going to a nearby park is also a fun thing to do!"

# Row 22 (ID 20)
$ws4.Cells.Item(22,1).Value = 20
$ws4.Cells.Item(22,2).Value = "no"
$ws4.Cells.Item(22,3).Value = "Message"
$ws4.Cells.Item(22,4).Value = "That's fine! You can always come back for more ideas if you need"

# Fill in the remaining ID / Type columns for rows 23-24
$ws4.Cells.Item(23,1).Value = 21
$ws4.Cells.Item(23,3).Value = "Message"
$ws4.Cells.Item(24,1).Value = 22
$ws4.Cells.Item(24,3).Value = "Message"

# ---------------------------------------------------------------------
# Selections / active sheet - the author finished up on Follow_Up_Messages
# ---------------------------------------------------------------------

$ws3.Range("F4").Select()
$ws4.Activate()
$ws4.Range("D27").Select()
